# 06.DanhSachChucNang.xlsx - "Cap nhat danh sach chuc nang"
#
# - Row 15, col G: reassign from "Huy" to "Nhi"
# - Rows 18,19,20,26,27: mark "Hoan thanh" (completion) column E as 100%
# - Row 47: mark "Hoan thanh" column E as 30%
# - Update the active view/selection (scrolled up, cell G14 selected)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Reassignment of function #12 ("Xem lich chay ke tiep") from Huy to Nhi
$ws.Range("G15").Value = "Nhi"

# Completion percentages (numFmt is already 0% on this column)
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("E47").Value = 0.3

# Restore the scrolled view / active selection as last left by the editor
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G14").Select()
